$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.018404666666666
$ws.Range("H2").Value = 12.055214
$ws.Range("I2").Value = 0.1784894308593523
$ws.Range("J2").Value = 0.1784894308593523
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.5736516666666667
$ws.Range("N2").Value = 1.720955
$ws.Range("O2").Value = 0.002679700009502673
$ws.Range("P2").Value = 0.002679700009502673
$ws.Range("Q2").Value = 2.305164534374444
$ws.Range("R2").Value = 20.74648080937
$ws.Range("S2").Value = 0.0004782981295699332
$ws.Range("T2").Value = 0.0004782981295699331

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.018404666666666
$ws.Range("H3").Value = 12.055214
$ws.Range("I3").Value = 0.1784894308593523
$ws.Range("J3").Value = 0.1784894308593523
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.928236333333333
$ws.Range("N3").Value = 11.784709
$ws.Range("O3").Value = 0.01834997708788796
$ws.Range("P3").Value = 0.01834997708788796
$ws.Range("Q3").Value = 15.78524321363622
$ws.Range("R3").Value = 142.067188922726
$ws.Range("S3").Value = 0.003275276966699278
$ws.Range("T3").Value = 0.003275276966699278

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.018404666666666
$ws.Range("H4").Value = 12.055214
$ws.Range("I4").Value = 0.1784894308593523
$ws.Range("J4").Value = 0.1784894308593523
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 205.3838603333334
$ws.Range("N4").Value = 616.1515810000001
$ws.Range("O4").Value = 0.9594099772863248
$ws.Range("P4").Value = 0.9594099772863247
$ws.Range("Q4").Value = 825.3154628214816
$ws.Range("R4").Value = 7427.839165393335
$ws.Range("S4").Value = 0.1712445408066203
$ws.Range("T4").Value = 0.1712445408066202

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.018404666666666
$ws.Range("H5").Value = 12.055214
$ws.Range("I5").Value = 0.1784894308593523
$ws.Range("J5").Value = 0.1784894308593523
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.187343666666667
$ws.Range("N5").Value = 12.562031
$ws.Range("O5").Value = 0.01956034561628449
$ws.Range("P5").Value = 0.01956034561628449
$ws.Range("Q5").Value = 16.82644133107044
$ws.Range("R5").Value = 151.437971979634
$ws.Range("S5").Value = 0.003491314956462847
$ws.Range("T5").Value = 0.003491314956462846

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.352037
$ws.Range("H6").Value = 13.056111
$ws.Range("I6").Value = 0.1933087062267439
$ws.Range("J6").Value = 0.1933087062267439
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.5736516666666667
$ws.Range("N6").Value = 1.720955
$ws.Range("O6").Value = 0.002679700009502673
$ws.Range("P6").Value = 0.002679700009502673
$ws.Range("Q6").Value = 2.496553278445
$ws.Range("R6").Value = 22.468979506005
$ws.Range("S6").Value = 0.000518009341912755
$ws.Range("T6").Value = 0.0005180093419127549

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.352037
$ws.Range("H7").Value = 13.056111
$ws.Range("I7").Value = 0.1933087062267439
$ws.Range("J7").Value = 0.1933087062267439
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.928236333333333
$ws.Range("N7").Value = 11.784709
$ws.Range("O7").Value = 0.01834997708788796
$ws.Range("P7").Value = 0.01834997708788796
$ws.Range("Q7").Value = 17.095829867411
$ws.Range("R7").Value = 153.862468806699
$ws.Range("S7").Value = 0.003547210330150015
$ws.Range("T7").Value = 0.003547210330150015

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.352037
$ws.Range("H8").Value = 13.056111
$ws.Range("I8").Value = 0.1933087062267439
$ws.Range("J8").Value = 0.1933087062267439
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 205.3838603333334
$ws.Range("N8").Value = 616.1515810000001
$ws.Range("O8").Value = 0.9594099772863248
$ws.Range("P8").Value = 0.9594099772863247
$ws.Range("Q8").Value = 893.8381593734991
$ws.Range("R8").Value = 8044.543434361492
$ws.Range("S8").Value = 0.1854623014502492
$ws.Range("T8").Value = 0.1854623014502491

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.352037
$ws.Range("H9").Value = 13.056111
$ws.Range("I9").Value = 0.1933087062267439
$ws.Range("J9").Value = 0.1933087062267439
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.187343666666667
$ws.Range("N9").Value = 12.562031
$ws.Range("O9").Value = 0.01956034561628449
$ws.Range("P9").Value = 0.01956034561628449
$ws.Range("Q9").Value = 18.223474569049
$ws.Range("R9").Value = 164.011271121441
$ws.Range("S9").Value = 0.003781185104431916
$ws.Range("T9").Value = 0.003781185104431916

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.14296233333333
$ws.Range("H10").Value = 42.428887
$ws.Range("I10").Value = 0.6282018629139038
$ws.Range("J10").Value = 0.6282018629139038
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.5736516666666667
$ws.Range("N10").Value = 1.720955
$ws.Range("O10").Value = 0.002679700009502673
$ws.Range("P10").Value = 0.002679700009502673
$ws.Range("Q10").Value = 8.113133914120555
$ws.Range("R10").Value = 73.01820522708499
$ws.Range("S10").Value = 0.001683392538019985
$ws.Range("T10").Value = 0.001683392538019985

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 14.14296233333333
$ws.Range("H11").Value = 42.428887
$ws.Range("I11").Value = 0.6282018629139038
$ws.Range("J11").Value = 0.6282018629139038
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.928236333333333
$ws.Range("N11").Value = 11.784709
$ws.Range("O11").Value = 0.01834997708788796
$ws.Range("P11").Value = 0.01834997708788796
$ws.Range("Q11").Value = 55.55689849876477
$ws.Range("R11").Value = 500.0120864888829
$ws.Range("S11").Value = 0.01152748979103867
$ws.Range("T11").Value = 0.01152748979103867

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 14.14296233333333
$ws.Range("H12").Value = 42.428887
$ws.Range("I12").Value = 0.6282018629139038
$ws.Range("J12").Value = 0.6282018629139038
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 205.3838603333334
$ws.Range("N12").Value = 616.1515810000001
$ws.Range("O12").Value = 0.9594099772863248
$ws.Range("P12").Value = 0.9594099772863247
$ws.Range("Q12").Value = 2904.736200568927
$ws.Range("R12").Value = 26142.62580512035
$ws.Range("S12").Value = 0.6027031350294554
$ws.Range("T12").Value = 0.6027031350294553

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 14.14296233333333
$ws.Range("H13").Value = 42.428887
$ws.Range("I13").Value = 0.6282018629139038
$ws.Range("J13").Value = 0.6282018629139038
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.187343666666667
$ws.Range("N13").Value = 12.562031
$ws.Range("O13").Value = 0.01956034561628449
$ws.Range("P13").Value = 0.01956034561628449
$ws.Range("Q13").Value = 59.22144375438855
$ws.Range("R13").Value = 532.992993789497
$ws.Range("S13").Value = 0.01228784555538973
$ws.Range("T13").Value = 0.01228784555538973
